# bi_generic_import sample product_variant.xlsx update (Nov 10)
# - Adds new "variant" related columns (Attribute / Variant Value / Variant price /
#   analytic_account_id / Expense Account / Income Account / Routes) in Z:AF
# - Tweaks a few BARCODE values in column E
# - Drops the stray trailing empty row
# - Minor view/row-height housekeeping

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New header row (Z1:AF1)
# ---------------------------------------------------------------------------
$ws.Range("Z1").Value = "Attribute"
$ws.Range("AA1").Value = "Variant Value"
$ws.Range("AB1").Value = "Variant price (Value Price Extra)"
$ws.Range("AC1").Value = "analytic_account_id"
$ws.Range("AD1").Value = "Expense Account (property_account_Expense_id)"
$ws.Range("AE1").Value = "Income Account (property_account_Income_id)"
$ws.Range("AF1").Value = "Routes (Route_ids)"

# ---------------------------------------------------------------------------
# 2. Row 2 (HP0012 / monitor - red) : barcode tweak + new variant data
# ---------------------------------------------------------------------------
$ws.Range("E2").Value = "34556f65"

$ws.Range("Z2").Value = "Color"
$ws.Range("AA2").Value = "White"
$ws.Range("AB2").Value = 100
$ws.Range("AC2").Value = "Administrative"
$ws.Range("AD2").Value = 101401
$ws.Range("AE2").Value = 110100
$ws.Range("AF2").Value = "Buy"

# ---------------------------------------------------------------------------
# 3. Row 3 (HP0013 / monitor - blue) : barcode tweak + new variant data
# ---------------------------------------------------------------------------
$ws.Range("E3").Value = "99677f615"

$ws.Range("Z3").Value = "Color"
$ws.Range("AA3").Value = "Black"
$ws.Range("AB3").Value = 150
$ws.Range("AC3").Value = "Administrative"
$ws.Range("AD3").Value = 101401
$ws.Range("AE3").Value = 110100
$ws.Range("AF3").Value = "Manufacture"

# ---------------------------------------------------------------------------
# 4. Row 4 (ST0011 / BOX - purple) : barcode tweak + new variant data
# ---------------------------------------------------------------------------
$ws.Range("E4").Value = "558585f5585"

$ws.Range("Z4").Value = "Legs"
$ws.Range("AA4").Value = "Steel"
$ws.Range("AB4").Value = 300
$ws.Range("AC4").Value = "Administrative"
$ws.Range("AD4").Value = 101401
$ws.Range("AE4").Value = 110100
$ws.Range("AF4").Value = "Manufacture"

# ---------------------------------------------------------------------------
# 5. Row 5 (PH0001 / Smart Phone - white) : barcode tweak + new variant data
# ---------------------------------------------------------------------------
$ws.Range("E5").Value = "B1110f0051"

$ws.Range("Z5").Value = "Legs"
$ws.Range("AA5").Value = "Custom"
$ws.Range("AB5").Value = 400
$ws.Range("AC5").Value = "Administrative"
$ws.Range("AD5").Value = 101401
$ws.Range("AE5").Value = 110100
$ws.Range("AF5").Value = "Manufacture"

# ---------------------------------------------------------------------------
# 6. Drop the stray trailing empty row 6 and normalise row heights back to
#    the sheet default (12.8) for the data rows.
# ---------------------------------------------------------------------------
$ws.Rows.Item(6).Delete()

$ws.Rows.Item(1).RowHeight = 12.8
$ws.Rows.Item(2).RowHeight = 12.8
$ws.Rows.Item(3).RowHeight = 12.8
$ws.Rows.Item(4).RowHeight = 12.8
$ws.Rows.Item(5).RowHeight = 12.8

# ---------------------------------------------------------------------------
# 7. Selection / view housekeeping to match the new, wider used range.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("I13").Select()
